# Add a new "robobet" entry to the deneme bonusu list.
# The source workbook keeps its rows alphabetically sorted (via AutoFilter +
# sortState on A:D), so the new record for "robobet" lands at row 332
# (between "risebet" and "roketbahis"), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 332 - shifts rows 332..410 down to 333..411.
$ws.Rows.Item(332).Insert() | Out-Null

# Fill in the new record: Site Adi | Sartlar | Durum | Yeni
$ws.Cells.Item(332, 1).Value = "robobet"
$ws.Cells.Item(332, 2).Value = "Maks 1.5k çekim"
$ws.Cells.Item(332, 3).Value = "yatırımsız"
$ws.Cells.Item(332, 4).Value = "Evet"

# Restore the active selection to A4 (matches the saved workbook UI state).
$ws.Range("A4").Select() | Out-Null
